$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 103: write A103 as text (matches existing date-label column A style) ---
$ws.Range("A103").NumberFormat = "@"
$ws.Range("A103").Value = "01-04-2021"
$ws.Range("A103").Style = $ws.Range("A102").Style

# --- Update existing row 102 values ---
$ws.Range("B102").Value = 152.8
$ws.Range("C102").Value = 85.90000000000001
$ws.Range("D102").Value = 223.3
$ws.Range("E102").Value = 228.2
$ws.Range("F102").Value = 180.7
$ws.Range("G102").Value = 129.7
$ws.Range("H102").Value = 145.7
$ws.Range("I102").Value = 146
$ws.Range("J102").Value = 144.4
$ws.Range("K102").Value = 64.5
$ws.Range("L102").Value = 173.4
$ws.Range("M102").Value = 133.2
$ws.Range("N102").Value = 113.6
$ws.Range("O102").Value = 108.1
$ws.Range("Q102").Value = 99.2
$ws.Range("R102").Value = 116.5
$ws.Range("S102").Value = 146.1
$ws.Range("T102").Value = 150.9
$ws.Range("U102").Value = 143
$ws.Range("V102").Value = 139.4
$ws.Range("W102").Value = 158.3
$ws.Range("X102").Value = 141.5
$ws.Range("Y102").Value = 80.2
$ws.Range("AA102").Value = 110.5
$ws.Range("AB102").Value = 138.3
$ws.Range("AC102").Value = 142.7
$ws.Range("AD102").Value = 187.6
$ws.Range("AE102").Value = 136.7
$ws.Range("AF102").Value = 135.1
$ws.Range("AG102").Value = 105.4
$ws.Range("AH102").Value = 148.2

# --- Add new row 103 values (B..AH) ---
$ws.Range("B103").Value = 161.9
$ws.Range("C103").Value = 74.2
$ws.Range("D103").Value = 230.7
$ws.Range("E103").Value = 234.5
$ws.Range("F103").Value = 199.9
$ws.Range("G103").Value = 138
$ws.Range("H103").Value = 152.7
$ws.Range("I103").Value = 156.6
$ws.Range("J103").Value = 145.3
$ws.Range("K103").Value = 68.5
$ws.Range("L103").Value = 174.8
$ws.Range("M103").Value = 150.8
$ws.Range("N103").Value = 129.6
$ws.Range("O103").Value = 147.4
$ws.Range("P103").Value = 121.5
$ws.Range("Q103").Value = 118.7
$ws.Range("R103").Value = 113.3
$ws.Range("S103").Value = 139
$ws.Range("T103").Value = 152.8
$ws.Range("U103").Value = 146.1
$ws.Range("V103").Value = 142.8
$ws.Range("W103").Value = 158.7
$ws.Range("X103").Value = 134
$ws.Range("Y103").Value = 74.2
$ws.Range("Z103").Value = 128.8
$ws.Range("AA103").Value = 111.2
$ws.Range("AB103").Value = 137.7
$ws.Range("AC103").Value = 142
$ws.Range("AD103").Value = 141
$ws.Range("AE103").Value = 133.1
$ws.Range("AF103").Value = 134.8
$ws.Range("AG103").Value = 122.4
$ws.Range("AH103").Value = 144.7

Write-Host "Update applied"
